# Auto-generated edit script
# Updates numeric cell values on multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# reflecting a recalculated 'Golem Profits' dataset refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 98
$ws.Range("H98").Value = 2226.75
$ws.Range("J98").Value = 2394
$ws.Range("L98").Value = 2394
$ws.Range("N98").Value = -5390

# Row 116
$ws.Range("H116").Value = 4500
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4500
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = $null
$ws.Range("N116").Value = -11384

# Row 122
$ws.Range("H122").Value = 2226.75
$ws.Range("J122").Value = 2394
$ws.Range("L122").Value = 7182
$ws.Range("N122").Value = -12082

# Row 125
$ws.Range("H125").Value = 1483
$ws.Range("I125").Value = 1251.4445
$ws.Range("J125").Value = 1899.8
$ws.Range("K125").Value = 11263.0005
$ws.Range("L125").Value = 17098.2
$ws.Range("M125").Value = -8803.0005
$ws.Range("N125").Value = -22018.2

# Row 137
$ws.Range("H137").Value = 2299.3333
$ws.Range("I137").Value = 1910.5
$ws.Range("J137").Value = 3077
$ws.Range("K137").Value = 5731.5
$ws.Range("L137").Value = 9231
$ws.Range("M137").Value = -3181.5
$ws.Range("N137").Value = -14331

# Row 138
$ws.Range("H138").Value = 7192.069
$ws.Range("J138").Value = 7214.84
$ws.Range("L138").Value = 21644.52
$ws.Range("N138").Value = -31924.52

# Row 140
$ws.Range("H140").Value = 95000
$ws.Range("J140").Value = 95000
$ws.Range("L140").Value = 95000
$ws.Range("N140").Value = -105360

$ws = $wb.Worksheets.Item("ARM")

# Row 74
$ws.Range("H74").Value = 4011.5
$ws.Range("I74").Value = 4011.5
$ws.Range("K74").Value = 4011.5
$ws.Range("M74").Value = -3137.5

# Row 77
$ws.Range("H77").Value = 4011.5
$ws.Range("I77").Value = 4011.5
$ws.Range("K77").Value = 20057.5
$ws.Range("M77").Value = -15689.5

# Row 132
$ws.Range("H132").Value = 3524.5
$ws.Range("I132").Value = 2699.3333
$ws.Range("K132").Value = 8097.999899999999
$ws.Range("M132").Value = -5567.999899999999

$ws = $wb.Worksheets.Item("BSM")

# Row 99
$ws.Range("H99").Value = 3773.1667
$ws.Range("I99").Value = 3888.9092
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 3888.9092
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -2390.9092
$ws.Range("N99").Value = -5496

# Row 112
$ws.Range("H112").Value = 195000
$ws.Range("J112").Value = 195000
$ws.Range("L112").Value = 195000
$ws.Range("N112").Value = -197954

# Row 141
$ws.Range("H141").Value = 20000
$ws.Range("I141").Value = 20000
$ws.Range("K141").Value = 20000
$ws.Range("M141").Value = -14820

$ws = $wb.Worksheets.Item("CRP")

# Row 8
$ws.Range("H8").Value = 6385.364
$ws.Range("I8").Value = 3096
$ws.Range("J8").Value = 9126.5
$ws.Range("K8").Value = 3096
$ws.Range("L8").Value = 9126.5
$ws.Range("M8").Value = -2956
$ws.Range("N8").Value = -9406.5

# Row 25
$ws.Range("H25").Value = 13005.5
$ws.Range("I25").Value = 1011
$ws.Range("K25").Value = 1011
$ws.Range("M25").Value = -837

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null

# Row 99
$ws.Range("H99").Value = 1002742.4
$ws.Range("J99").Value = 1856
$ws.Range("L99").Value = 1856
$ws.Range("N99").Value = -4852

# Row 126
$ws.Range("H126").Value = 1002742.4
$ws.Range("J126").Value = 1856
$ws.Range("L126").Value = 5568
$ws.Range("N126").Value = -10508

# Row 132
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

# Row 141
$ws.Range("H141").Value = 832665.8
$ws.Range("J141").Value = 832665.8
$ws.Range("L141").Value = 832665.8
$ws.Range("N141").Value = -843025.8

$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 26740.55
$ws.Range("I4").Value = 46164.637
$ws.Range("K4").Value = 138493.911
$ws.Range("M4").Value = -138381.911

# Row 44
$ws.Range("H44").Value = 282.5
$ws.Range("I44").Value = 282.5
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 847.5
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -449.5
$ws.Range("N44").Value = $null

# Row 129
$ws.Range("H129").Value = 2627.4
$ws.Range("I129").Value = 2739.125
$ws.Range("K129").Value = 8217.375
$ws.Range("M129").Value = -3217.375

# Row 137
$ws.Range("H137").Value = 3999.5
$ws.Range("I137").Value = 3999.5
$ws.Range("K137").Value = 11998.5
$ws.Range("M137").Value = -6898.5

$ws = $wb.Worksheets.Item("GSM")

# Row 2
$ws.Range("H2").Value = 130.88235
$ws.Range("I2").Value = 110.416664
$ws.Range("K2").Value = 110.416664
$ws.Range("M2").Value = 2.583336000000003

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = $null

# Row 43
$ws.Range("H43").Value = 6373
$ws.Range("I43").Value = 6373
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 6373
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -6222
$ws.Range("N43").Value = $null

# Row 46
$ws.Range("H46").Value = 12000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null

# Row 57
$ws.Range("H57").Value = 9590
$ws.Range("J57").Value = 35000
$ws.Range("L57").Value = 35000
$ws.Range("N57").Value = -36640

# Row 126
$ws.Range("H126").Value = 1807
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2114
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 6342
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -11282

$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4888

# Row 122
$ws.Range("H122").Value = 5499.75
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 17998.5
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -22898.5

# Row 126
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 2139
$ws.Range("I122").Value = 1933.3334
$ws.Range("J122").Value = 2447.5
$ws.Range("K122").Value = 5800.0002
$ws.Range("L122").Value = 7342.5
$ws.Range("M122").Value = -3350.0002
$ws.Range("N122").Value = -12242.5

# Row 126
$ws.Range("H126").Value = 2806.7144
$ws.Range("I126").Value = 1629.4
$ws.Range("K126").Value = 4888.200000000001
$ws.Range("M126").Value = -2418.200000000001
